# Update profit files after running on 2025-10-20
# Adds a new data row (row 49) to Sheet1 with the latest BTC/KAS allocation split.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a leading apostrophe so Excel stores the date as literal text
# (matching the existing Date column cells) instead of auto-converting
# it to a date serial number; resetting the style afterwards keeps the
# new row's formatting identical to the other un-styled data rows.
$ws.Range("A49").Value = "'10/20/2025"
$ws.Range("A49").Style = "Normal"
$ws.Range("B49").Value = 0.1910841114775904
$ws.Range("C49").Value = 0.8089158885224096
